$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 34, shifting existing rows 34-92 down to 35-93.
$ws.Rows.Item(34).Insert()

# Populate the new row 34 with the new data point.
$ws.Cells.Item(34, 1).Value = 9
$ws.Cells.Item(34, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(34, 3).Value = "Metropolitana"
$ws.Cells.Item(34, 4).Value = 45100
$ws.Cells.Item(34, 5).Value = 13
$ws.Cells.Item(34, 6).Value = 100112035
$ws.Cells.Item(34, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 43
$ws.Cells.Item(34, 11).Value = 17000
$ws.Cells.Item(34, 12).Value = 18000
$ws.Cells.Item(34, 13).Value = 17512
$ws.Cells.Item(34, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(34, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(34, 16).Value = 1167
$ws.Cells.Item(34, 17).Value = 15
$ws.Cells.Item(34, 18).Value = "Hortaliza"
